# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.8
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 12
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 67

# Row 6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11

# Row 11
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 1.33
$ws.Range("J11").Value = 6.3
$ws.Range("K11").Value = 2.7
$ws.Range("P11").Value = 5.7
$ws.Range("Q11").Value = 1.35
$ws.Range("R11").Value = 2.95
$ws.Range("S11").Value = 1.21
$ws.Range("T11").Value = 3.95
$ws.Range("V11").Value = 2.27
$ws.Range("W11").Value = 35
$ws.Range("AA11").Value = 70
$ws.Range("AB11").Value = 45
$ws.Range("AC11").Value = 10.25
$ws.Range("AD11").Value = 11
$ws.Range("AG11").Value = 11.5
$ws.Range("AH11").Value = 9
$ws.Range("AI11").Value = 8.75
$ws.Range("AL11").Value = 17.5
$ws.Range("AM11").Value = 200
$ws.Range("AO11").Value = 37
$ws.Range("AP11").Value = 29
$ws.Range("AQ11").Value = 200
$ws.Range("AT11").Value = 3.95
$ws.Range("AV11").Value = 40

$wb.Save()
